$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4").Value = -7.904600000000007
$ws.Range("C9").Value = -13.4575
$ws.Range("D9").Value = -8.289500000000002
$ws.Range("D11").Value = -8.520699999999994
$ws.Range("C18").Value = -12.01449999999998
$ws.Range("C20").Value = -11.3878
$ws.Range("D23").Value = -7.9908
$ws.Range("D24").Value = -7.169700000000001
$ws.Range("D26").Value = -7.557200000000002
$ws.Range("C27").Value = -12.0777
$ws.Range("D34").Value = -8.111600000000001
$ws.Range("C35").Value = -11.872
$ws.Range("D35").Value = -8.311199999999996
$ws.Range("D48").Value = -8.126300000000002
$ws.Range("D49").Value = -7.8411
$ws.Range("D52").Value = -7.992800000000005
$ws.Range("D66").Value = -7.018899999999999
$ws.Range("D67").Value = -7.426399999999994
$ws.Range("C69").Value = -10.49249999999999
$ws.Range("C76").Value = -12.59020000000001
$ws.Range("C78").Value = -11.94399999999999
$ws.Range("D78").Value = -8.126200000000004
$ws.Range("D80").Value = -7.753600000000001
$ws.Range("C82").Value = -11.9884
$ws.Range("C83").Value = -13.2163
$ws.Range("C93").Value = -10.1701
$ws.Range("D99").Value = -8.137500000000006
$ws.Range("D104").Value = -7.702600000000003
